# Natmi following Dr Hou advice
# Update Il12b-Il12rb2 LR-pair table: recompute edge weights for the M1/M2 sending
# clusters against all six target clusters (ECs, FAPs, M1, M2, Neutro, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(@("M1", "Il12b", "Il12rb2", "ECs"), @(1, 0.3333333333333333, 0.2137916666666667, 0.641375, 0.1880241480025247, 0.1880241480025247, 2, 1, 0.2521135, 0.504227, 0.05590272485875765, 0.0405525532384611, 0.05389976535416667, 0.323398592125, 0.01051106221258746, 0.007624859271988671)),
    @(@("M1", "Il12b", "Il12rb2", "FAPs"), @(1, 0.3333333333333333, 0.2137916666666667, 0.641375, 0.1880241480025247, 0.1880241480025247, 3, 1, 0.5416043333333334, 1.624813, 0.1200933628252521, 0.1306758973340255, 0.1157904930972222, 1.042114437875, 0.0225804522259761, 0.02457022426069554)),
    @(@("M1", "Il12b", "Il12rb2", "M1"), @(1, 0.3333333333333333, 0.2137916666666667, 0.641375, 0.1880241480025247, 0.1880241480025247, 3, 1, 1.132338333333333, 3.397015, 0.2510805581428902, 0.2732055832776724, 0.2420844995138889, 2.178760495625, 0.04720920802481528, 0.05136924702531716)),
    @(@("M1", "Il12b", "Il12rb2", "M2"), @(1, 0.3333333333333333, 0.2137916666666667, 0.641375, 0.1880241480025247, 0.1880241480025247, 3, 1, 1.434601, 4.303803, 0.318103175693085, 0.3461341821944255, 0.3067057387916667, 2.760351649125, 0.05981107858658972, 0.06508158470165751)),
    @(@("M1", "Il12b", "Il12rb2", "Neutro"), @(1, 0.3333333333333333, 0.2137916666666667, 0.641375, 0.1880241480025247, 0.1880241480025247, 3, 1, 0.30565, 0.9169499999999999, 0.06777371244728772, 0.07374587971688724, 0.06534542291666667, 0.58810880625, 0.01274309453986937, 0.01386600620246439)),
    @(@("M1", "Il12b", "Il12rb2", "sCs"), @(1, 0.3333333333333333, 0.2137916666666667, 0.641375, 0.1880241480025247, 0.1880241480025247, 2, 1, 0.8435535000000001, 1.687107, 0.1870464660327274, 0.1356859042385283, 0.1803447086875, 1.082068252125, 0.03516925241268674, 0.02551222654040143)),
    @(@("M2", "Il12b", "Il12rb2", "ECs"), @(3, 1, 0.9232520000000001, 2.769756, 0.8119758519974754, 0.8119758519974754, 2, 1, 0.2521135, 0.504227, 0.05590272485875765, 0.0405525532384611, 0.232764293102, 1.396585758612, 0.04539166264617019, 0.03292769396647242)),
    @(@("M2", "Il12b", "Il12rb2", "FAPs"), @(3, 1, 0.9232520000000001, 2.769756, 0.8119758519974754, 0.8119758519974754, 3, 1, 0.5416043333333334, 1.624813, 0.1200933628252521, 0.1306758973340255, 0.5000372839586668, 4.500335555628, 0.09751291059927601, 0.10610567307333)),
    @(@("M2", "Il12b", "Il12rb2", "M1"), @(3, 1, 0.9232520000000001, 2.769756, 0.8119758519974754, 0.8119758519974754, 3, 1, 1.132338333333333, 3.397015, 0.2510805581428902, 0.2732055832776724, 1.045433630926667, 9.408902678339999, 0.2038713501180749, 0.2218363362523553)),
    @(@("M2", "Il12b", "Il12rb2", "M2"), @(3, 1, 0.9232520000000001, 2.769756, 0.8119758519974754, 0.8119758519974754, 3, 1, 1.434601, 4.303803, 0.318103175693085, 0.3461341821944255, 1.324498242452, 11.920484182068, 0.2582920971064953, 0.2810525974927681)),
    @(@("M2", "Il12b", "Il12rb2", "Neutro"), @(3, 1, 0.9232520000000001, 2.769756, 0.8119758519974754, 0.8119758519974754, 3, 1, 0.30565, 0.9169499999999999, 0.06777371244728772, 0.07374587971688724, 0.2821919738, 2.5397277642, 0.05503061790741835, 0.05987987351442286)),
    @(@("M2", "Il12b", "Il12rb2", "sCs"), @(3, 1, 0.9232520000000001, 2.769756, 0.8119758519974754, 0.8119758519974754, 2, 1, 0.8435535000000001, 1.687107, 0.1870464660327274, 0.1356859042385283, 0.7788124559820001, 4.672874735892001, 0.1518772136200407, 0.1101736776981268))
)

$r = 2
foreach ($row in $rows) {
    $strs = $row[0]
    $vals = $row[1]
    for ($i = 0; $i -lt $strs.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $strs[$i]
    }
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($r, $strs.Count + $j + 1).Value = $vals[$j]
    }
    $r++
}
